$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 300
$ws.Range("I32").Value = 300
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 300
$ws.Range("L32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = 26

$ws.Range("H33").Value = 123.92857
$ws.Range("I33").Value = 87.38461
$ws.Range("K33").Value = 87.38461
$ws.Range("M33").Value = 141.61539

$ws.Range("H40").Value = 6210.931
$ws.Range("J40").Value = 6820.294
$ws.Range("L40").Value = 6820.294
$ws.Range("N40").Value = -7170.294

$ws.Range("H58").Value = 1492.75
$ws.Range("I58").Value = 203.5
$ws.Range("J58").Value = 2137.375
$ws.Range("K58").Value = 610.5
$ws.Range("L58").Value = 6412.125
$ws.Range("M58").Value = -460.5
$ws.Range("N58").Value = -6712.125

$ws.Range("H62").Value = 11200
$ws.Range("I62").Value = 8000
$ws.Range("K62").Value = 8000
$ws.Range("M62").Value = -7376

$ws.Range("H65").Value = 11200
$ws.Range("I65").Value = 8000
$ws.Range("K65").Value = 40000
$ws.Range("M65").Value = -36880

$ws.Range("H70").Value = 6707.2856
$ws.Range("I70").Value = 3580.6
$ws.Range("K70").Value = 10741.8
$ws.Range("M70").Value = -10471.8

$ws.Range("H73").Value = 6707.2856
$ws.Range("I73").Value = 3580.6
$ws.Range("K73").Value = 10741.8
$ws.Range("M73").Value = -9805.799999999999

$ws.Range("H86").Value = 3656.6
$ws.Range("J86").Value = 3965.6667
$ws.Range("L86").Value = 3965.6667
$ws.Range("N86").Value = -6211.6667

$ws.Range("H89").Value = 3656.6
$ws.Range("J89").Value = 3965.6667
$ws.Range("L89").Value = 19828.3335
$ws.Range("N89").Value = -31060.3335

$ws.Range("H116").Value = 6054.6665
$ws.Range("I116").Value = 5213
$ws.Range("J116").Value = 9000.5
$ws.Range("K116").Value = 5213
$ws.Range("L116").Value = 9000.5
$ws.Range("M116").Value = -1771
$ws.Range("N116").Value = -15884.5

$ws.Range("H127").Value = 9966.666999999999
$ws.Range("J127").Value = 9950
$ws.Range("L127").Value = 29850
$ws.Range("N127").Value = -39770

$ws.Range("H138").Value = 5863
$ws.Range("J138").Value = 6292.3335
$ws.Range("L138").Value = 18877.0005
$ws.Range("N138").Value = -29157.0005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1127.1428
$ws.Range("I2").Value = 1131.8334
$ws.Range("K2").Value = 1131.8334
$ws.Range("M2").Value = -1018.8334

$ws.Range("H45").Value = 2834.0625
$ws.Range("I45").Value = 2120.6365
$ws.Range("K45").Value = 2120.6365
$ws.Range("M45").Value = -1743.6365

$ws.Range("H116").Value = 1127.1428
$ws.Range("I116").Value = 1131.8334
$ws.Range("K116").Value = 1131.8334
$ws.Range("M116").Value = 1162.1666

$ws.Range("H132").Value = 2389
$ws.Range("I132").Value = 2401.6875
$ws.Range("K132").Value = 7205.0625
$ws.Range("M132").Value = -4675.0625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1127.1428
$ws.Range("I3").Value = 1131.8334
$ws.Range("K3").Value = 1131.8334
$ws.Range("M3").Value = -1017.8334

$ws.Range("H80").Value = 607
$ws.Range("I80").Value = 357.6
$ws.Range("K80").Value = 357.6
$ws.Range("M80").Value = 640.4

$ws.Range("H83").Value = 607
$ws.Range("I83").Value = 357.6
$ws.Range("K83").Value = 1788
$ws.Range("M83").Value = 3204

$ws.Range("H134").Value = 3916.3635
$ws.Range("I134").Value = 3898.889
$ws.Range("K134").Value = 11696.667
$ws.Range("M134").Value = -9161.667000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 774.25
$ws.Range("I16").Value = 699
$ws.Range("J16").Value = 1000
$ws.Range("K16").Value = 699
$ws.Range("L16").Value = 1000
$ws.Range("M16").Value = -412
$ws.Range("N16").Value = -1574

$ws.Range("H105").Value = 1690.5454
$ws.Range("I105").Value = 727.25
$ws.Range("J105").Value = 2241
$ws.Range("K105").Value = 727.25
$ws.Range("L105").Value = 2241
$ws.Range("M105").Value = 1019.75
$ws.Range("N105").Value = -5735

$ws.Range("H113").Value = 774.25
$ws.Range("I113").Value = 699
$ws.Range("J113").Value = 1000
$ws.Range("K113").Value = 699
$ws.Range("L113").Value = 1000
$ws.Range("M113").Value = 1471
$ws.Range("N113").Value = -5340

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 946.1
$ws.Range("J34").Value = 1483.1666
$ws.Range("L34").Value = 4449.4998
$ws.Range("N34").Value = -4617.4998

$ws.Range("H39").Value = 5026.7334
$ws.Range("I39").Value = 310
$ws.Range("J39").Value = 6205.9165
$ws.Range("K39").Value = 930
$ws.Range("L39").Value = 18617.7495
$ws.Range("M39").Value = -636
$ws.Range("N39").Value = -19205.7495

$ws.Range("H47").Value = 775.75
$ws.Range("I47").Value = 775.75
$ws.Range("K47").Value = 2327.25
$ws.Range("M47").Value = -1896.25

$ws.Range("H60").Value = 1490.9375
$ws.Range("I60").Value = 177
$ws.Range("J60").Value = 2279.3
$ws.Range("K60").Value = 531
$ws.Range("L60").Value = 6837.900000000001
$ws.Range("M60").Value = -280
$ws.Range("N60").Value = -7339.900000000001

$ws.Range("H68").Value = 3500
$ws.Range("I68").Value = 1998
$ws.Range("J68").Value = 3875.5
$ws.Range("K68").Value = 5994
$ws.Range("L68").Value = 11626.5
$ws.Range("M68").Value = -5183
$ws.Range("N68").Value = -13248.5

$ws.Range("H71").Value = 3500
$ws.Range("I71").Value = 1998
$ws.Range("J71").Value = 3875.5
$ws.Range("K71").Value = 17982
$ws.Range("L71").Value = 34879.5
$ws.Range("M71").Value = -13926
$ws.Range("N71").Value = -42991.5

$ws.Range("H113").Value = 296
$ws.Range("I113").Value = 400
$ws.Range("J113").Value = 261.33334
$ws.Range("K113").Value = 1200
$ws.Range("L113").Value = 784.0000200000001
$ws.Range("M113").Value = 970
$ws.Range("N113").Value = -5124.00002

$ws.Range("H114").Value = 4895
$ws.Range("J114").Value = 14031
$ws.Range("L114").Value = 42093
$ws.Range("N114").Value = -48601

$ws.Range("H136").Value = 7558.3335
$ws.Range("I136").Value = 6975
$ws.Range("K136").Value = 20925
$ws.Range("M136").Value = -15825

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3721.1052
$ws.Range("I61").Value = 2713.4
$ws.Range("K61").Value = 2713.4
$ws.Range("M61").Value = -2511.4

$ws.Range("H113").Value = 3721.1052
$ws.Range("I113").Value = 2713.4
$ws.Range("K113").Value = 2713.4
$ws.Range("M113").Value = -543.4000000000001

$ws.Range("H136").Value = 2078.6428
$ws.Range("I136").Value = 2078.6428
$ws.Range("K136").Value = 6235.928400000001
$ws.Range("M136").Value = -3685.928400000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 2492.8
$ws.Range("I107").Value = 2492.8
$ws.Range("K107").Value = 7478.400000000001
$ws.Range("M107").Value = -5558.400000000001

$ws.Range("H113").Value = 506.8889
$ws.Range("I113").Value = 392.8
$ws.Range("J113").Value = 649.5
$ws.Range("K113").Value = 1178.4
$ws.Range("L113").Value = 1948.5
$ws.Range("M113").Value = 991.5999999999999
$ws.Range("N113").Value = -6288.5
